$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the K column (Fecha) number format so dates render as
# "YYYY-MM-DD HH:MM:SS" (matches new numFmt 165), and write the real
# date value (2024-08-09) instead of the inline text "09/08/24".
$range = $ws.Range("K2:K23")
$range.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$range.Value = (Get-Date -Year 2024 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
